# Scrum_Log_Book.xlsx update
# - Fills in the previously-empty "E" column (developer activity) entries
#   on the ActivityLog sheet for rows 4-30.
# - Restores the active-cell selections on all three sheets.

$wb = $excel.ActiveWorkbook

$activityLog = $wb.Worksheets.Item("ActivityLog")

# Row-by-row E column updates (row -> new text)
$activityLog.Range("E4").Value  = "Share examples of JSON and XML formats found online with team."
$activityLog.Range("E5").Value  = "No Activity"
$activityLog.Range("E6").Value  = "Review of Google App Engine Code on Showing Taxi Location"
$activityLog.Range("E7").Value  = "Scrum meeting to share progress and to decide what to be done to complete Sprint 1."
$activityLog.Range("E8").Value  = "Scrum meeting to share progress and to decide what to be done to complete Sprint 1."
$activityLog.Range("E9").Value  = "No Activity"
$activityLog.Range("E10").Value = "Scrum meeting to share progress and to decide what to be done to complete Sprint 1. Research on displaying location on Google Map based on longitude and latitude."
$activityLog.Range("E11").Value = "Created Git Account to access and share progress."
$activityLog.Range("E12").Value = "Managed to implement showing taxi location on a Second map by looping through the taxi location from the JSON file"
$activityLog.Range("E13").Value = "No Activity"
$activityLog.Range("E14").Value = "Remote session with Team to discuss about Git's features and usage."
$activityLog.Range("E15").Value = "No Activity"
$activityLog.Range("E16").Value = "No Activity"
$activityLog.Range("E19").Value = "Working on code to present both Taxi Locations and Current User location to a single map. "
$activityLog.Range("E20").Value = "Successfully integrated both Taxi Locations and Current User location to a single map. Uploaded code to Git for sharing among team members. "
$activityLog.Range("E22").Value = "No Activity"
$activityLog.Range("E23").Value = "No Activity"
$activityLog.Range("E25").Value = "New requirement to change icon of taxi. Research and review of code for implementation. "
$activityLog.Range("E26").Value = "Managed to changed icon of all taxi location by retreiving image from web in PNG format to differentiate between current user's and taxi's."
$activityLog.Range("E27").Value = "Review UI design and submitted updated code to Git."
$activityLog.Range("E28").Value = "Tried to add a distance calculator API from google map to calculate distance between user and all taxis but was not successful."
$activityLog.Range("E29").Value = "No Activity"
$activityLog.Range("E30").Value = "Updated SCRUM master on the project progress."

# Restore the per-sheet active-cell selections. Select the non-active
# sheets first so that ActivityLog ends up as the final activated /
# tab-selected sheet (matching the original workbook's active tab).
$productBacklog = $wb.Worksheets.Item("ProductBacklog")
$productBacklog.Range("D19").Select() | Out-Null

$sprintBacklog = $wb.Worksheets.Item("SprintBacklog")
$sprintBacklog.Range("D25").Select() | Out-Null

$activityLog.Range("E4").Select() | Out-Null
